$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values could be misread as numbers by Excel;
# force them to Text format before assignment, then restore the
# original (Normal) style so no stray formatting is introduced.
$textCells = @("D4", "D5", "D7", "D9", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.784.61"
$ws.Range("E2").Value = "  -1.70%  "
$ws.Range("D3").Value = "1.871.09"
$ws.Range("E3").Value = "  -1.90%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "300.97"
$ws.Range("E5").Value = "  -2.17%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "0.5321"
$ws.Range("E7").Value = "  +1.29%  "
$ws.Range("E8").Value = "  -1.99%  "
$ws.Range("D9").Value = "0.07176"
$ws.Range("E9").Value = "  -1.81%  "
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("D11").Value = "0.8879"
$ws.Range("E11").Value = "  -1.97%  "
$ws.Range("D12").Value = "0.08161"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "1.884.79"
$ws.Range("E13").Value = "  +14.04%  "
$ws.Range("D14").Value = "92.77"
$ws.Range("E14").Value = "  -3.79%  "
$ws.Range("D15").Value = "5.287"
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "14.80"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "0.000008489"
$ws.Range("E18").Value = "  -2.31%  "
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "26.847.61"
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("D21").Value = "4.977"
$ws.Range("E21").Value = "  -2.89%  "
$ws.Range("D23").Value = "6.386"
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("D24").Value = "2.288"
$ws.Range("E24").Value = "  -2.59%  "
$ws.Range("D25").Value = "146.44"
$ws.Range("E25").Value = "  -2.34%  "
$ws.Range("D26").Value = "1.736"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").Value = "18.01"
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("D28").Value = "114.01"
$ws.Range("E28").Value = "  -2.46%  "
$ws.Range("D29").Value = "4.706"
$ws.Range("E29").Value = "  -2.91%  "
$ws.Range("D30").Value = "4.604"
$ws.Range("D31").Value = "0.09125"
$ws.Range("E31").Value = "  -1.38%  "
$ws.Range("D32").Value = "0.8082"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("D33").Value = "0.05010"
$ws.Range("E33").Value = "  -1.36%  "
$ws.Range("D34").Value = "1.167"
$ws.Range("E34").Value = "  -5.04%  "
$ws.Range("D35").Value = "2.960"
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("D36").Value = "0.6049"
$ws.Range("E36").Value = "  +5.01%  "
$ws.Range("D37").Value = "2.651"
$ws.Range("E37").Value = "  -2.71%  "
$ws.Range("E38").Value = "  -4.65%  "
$ws.Range("E39").Value = "  -2.48%  "
$ws.Range("D40").Value = "1.068"
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("D41").Value = "6.532"
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("D42").Value = "8.764"
$ws.Range("E42").Value = "  -3.36%  "
$ws.Range("D43").Value = "0.5144"
$ws.Range("E43").Value = "  +4.65%  "
$ws.Range("D44").Value = "115.02"
$ws.Range("E44").Value = "  -2.31%  "
$ws.Range("D45").Value = "0.1492"
$ws.Range("E45").Value = "  -2.00%  "
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("D48").Value = "9.925"
$ws.Range("E48").Value = "  -2.92%  "
$ws.Range("D49").Value = "37.45"
$ws.Range("E49").Value = "  -3.43%  "
$ws.Range("D50").Value = "0.06041"
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("D51").Value = "62.23"
$ws.Range("E51").Value = "  -3.68%  "

foreach ($cell in $textCells) {
    $ws.Range($cell).Style = "Normal"
}
